$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive component forecaster bug fix shifts the existing error rows down by
# one (each row's B:G values move into the next row), making room at row 2 for
# a freshly computed error row. Column A (the Q-period labels) stays fixed.
# Process bottom-up so source values aren't clobbered before being read.
for ($r = 10; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("B$dst").Value2 = $ws.Range("B$src").Value2
    $ws.Range("C$dst").Value2 = $ws.Range("C$src").Value2
    $ws.Range("D$dst").Value2 = $ws.Range("D$src").Value2
    $ws.Range("E$dst").Value2 = $ws.Range("E$src").Value2
    $ws.Range("F$dst").Value2 = $ws.Range("F$src").Value2
    $ws.Range("G$dst").Value2 = $ws.Range("G$src").Value2
}

# Write the newly computed forecaster error values into row 2.
$ws.Range("B2").Value2 = [double]"2.026185291058083E-07"
$ws.Range("C2").Value2 = [double]"6.752006303172386E-07"
$ws.Range("D2").Value2 = [double]"2.666983691600854E-12"
$ws.Range("E2").Value2 = [double]"1.6330902276362E-06"
$ws.Range("F2").Value2 = [double]"1.664876689494737E-06"
$ws.Range("G2").Value2 = 19
